$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17, shifting existing rows 17-101 down to 18-102
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new record
$ws.Range("A17").Value = 4
$ws.Range("B17").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C17").Value = "Los Lagos"
$ws.Range("D17").Value = 44819
$ws.Range("E17").Value = 10
$ws.Range("F17").Value = 100112026
$ws.Range("G17").Value = "Haba"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 80
$ws.Range("K17").Value = 14000
$ws.Range("L17").Value = 14000
$ws.Range("M17").Value = 14000
$ws.Range("N17").Value = "`$/saco 25 kilos"
$ws.Range("O17").Value = "Provincia de Limarí"
$ws.Range("P17").Value = 560
$ws.Range("Q17").Value = 25
$ws.Range("R17").Value = "Hortaliza"
